# Update decimals performance analysis
$wb = $excel.ActiveWorkbook

$wsResults = $wb.Worksheets.Item("results")
$wsScale2  = $wb.Worksheets.Item("scale=2")
$wsScale4  = $wb.Worksheets.Item("scale=4")
$wsScale6  = $wb.Worksheets.Item("scale=6")

# ---------------------------------------------------------------
# 1) Refresh the measured timings on the "results" sheet (column C)
# ---------------------------------------------------------------
$wsResults.Range("C2").Value = 1.834
$wsResults.Range("C3").Value = 1.531
$wsResults.Range("C4").Value = 1.746
$wsResults.Range("C5").Value = 1.51
$wsResults.Range("C6").Value = 1.661
$wsResults.Range("C7").Value = 1.533
$wsResults.Range("C8").Value = 1.696
$wsResults.Range("C9").Value = 1.539
$wsResults.Range("C10").Value = 1.993
$wsResults.Range("C11").Value = 1.53
$wsResults.Range("C12").Value = 1.893
$wsResults.Range("C13").Value = 1.654
$wsResults.Range("C14").Value = 1.628
$wsResults.Range("C15").Value = 1.552
$wsResults.Range("C16").Value = 2.025
$wsResults.Range("C17").Value = 1.582
$wsResults.Range("C18").Value = 1.901
$wsResults.Range("C19").Value = 2.028
$wsResults.Range("C20").Value = 1.507
$wsResults.Range("C21").Value = 1.771
$wsResults.Range("C22").Value = 2.115
$wsResults.Range("C23").Value = 1.521
$wsResults.Range("C24").Value = 2.09
$wsResults.Range("C25").Value = 2.287
$wsResults.Range("C26").Value = 2.57
$wsResults.Range("C27").Value = 2.734
$wsResults.Range("C28").Value = 3.044
$wsResults.Range("C29").Value = 2.968
$wsResults.Range("C30").Value = 2.773
$wsResults.Range("C31").Value = 3.069
$wsResults.Range("C32").Value = 2.837
$wsResults.Range("C33").Value = 2.768
$wsResults.Range("C34").Value = 3.074
$wsResults.Range("C35").Value = 3.148
$wsResults.Range("C36").Value = 2.758
$wsResults.Range("C37").Value = 3.17
$wsResults.Range("C38").Value = 2.73
$wsResults.Range("C39").Value = 3.14
$wsResults.Range("C40").Value = 3.045

# Filter the "results" table down to scale == 6 (hides the other rows
# and writes the corresponding <autoFilter> definition).
$wsResults.Range("A1:C40").AutoFilter(2, @("6"), 7)

# ---------------------------------------------------------------
# 2) Mirror the same timing refresh on the per-scale breakout sheets
# ---------------------------------------------------------------
$wsScale2.Range("C2").Value = 1.834
$wsScale2.Range("C3").Value = 1.531
$wsScale2.Range("C4").Value = 1.51
$wsScale2.Range("C5").Value = 1.533
$wsScale2.Range("C6").Value = 1.539
$wsScale2.Range("C7").Value = 1.53
$wsScale2.Range("C8").Value = 1.654
$wsScale2.Range("C9").Value = 1.552
$wsScale2.Range("C10").Value = 1.582
$wsScale2.Range("C11").Value = 1.507
$wsScale2.Range("C12").Value = 1.521
$wsScale2.Range("C13").Value = 2.57
$wsScale2.Range("C14").Value = 2.968
$wsScale2.Range("C15").Value = 2.837

$wsScale4.Range("C2").Value = 1.746
$wsScale4.Range("C3").Value = 1.661
$wsScale4.Range("C4").Value = 1.696
$wsScale4.Range("C5").Value = 1.993
$wsScale4.Range("C6").Value = 1.893
$wsScale4.Range("C7").Value = 1.628
$wsScale4.Range("C8").Value = 2.025
$wsScale4.Range("C9").Value = 1.901
$wsScale4.Range("C10").Value = 1.771
$wsScale4.Range("C11").Value = 2.09
$wsScale4.Range("C12").Value = 2.734
$wsScale4.Range("C13").Value = 2.773
$wsScale4.Range("C14").Value = 2.768
$wsScale4.Range("C15").Value = 3.148
$wsScale4.Range("C16").Value = 3.17
$wsScale4.Range("C17").Value = 3.14

$wsScale6.Range("C2").Value = 2.028
$wsScale6.Range("C3").Value = 2.115
$wsScale6.Range("C4").Value = 2.287
$wsScale6.Range("C5").Value = 3.044
$wsScale6.Range("C6").Value = 3.069
$wsScale6.Range("C7").Value = 3.074
$wsScale6.Range("C8").Value = 2.758
$wsScale6.Range("C9").Value = 2.73
$wsScale6.Range("C10").Value = 3.045

# ---------------------------------------------------------------
# 3) Add a new "Sheet4" worksheet with a byte-count -> base-10 digit
#    precision table, placed after "scale=6".
# ---------------------------------------------------------------
$last  = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsNew.Name = "Sheet4"

$wsNew.Range("A1").Value = "bytes"
$wsNew.Range("B1").Value = "base-10 digits"

for ($i = 1; $i -le 16; $i++) {
    $row = $i + 1
    $wsNew.Range("A$row").Value = $i
    $wsNew.Range("B$row").Formula = "=FLOOR(LOG10(2^(8*A$row-1)-1),1)"
}

$wb.Save()
